$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the last data row (row 7): the aircraft code "RF12" is replaced
# by a new "RF32" entry whose full/empty tank capacities are both 600.
$ws.Range("A7").Value = "RF32"
$ws.Range("B7").Value = 600
$ws.Range("C7").Value = 600
